{"js": "// 1) Remove the \"Previous Participant Interest\" bullet paragraph entirely\n//    (it sat right after the \"Perigean Contract\" bullet, inside the\n//    \"What is your recruitment strategy?\" list).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText =\n  \"Previous Participant Interest \\u2013 We have conducted user feedback \" +\n  \"sessions in the past for DGIB and have participants who have \" +\n  \"indicated that they would like to participate in future research. \";\n\nlet removed = false;\nlet teamRolesParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (!removed && para.text === targetText) {\n    para.delete();\n    removed = true;\n    continue;\n  }\n  if (para.text === \"Team Roles\") {\n    teamRolesParagraph = para;\n  }\n}\nawait context.sync();\n\nif (!removed) {\n  throw new Error('Could not find the \"Previous Participant Interest\" paragraph to remove.');\n}\nif (!teamRolesParagraph) {\n  throw new Error('Could not find the \"Team Roles\" heading paragraph.');\n}\n\n// 2) Mark the \"Team Roles\" heading run with a lastRenderedPageBreak\n//    (Word stamps this on the run that happens to start a new page;\n//    it is a pure rendering hint with no visible effect on the text).\n//    We rebuild the paragraph via insertOoxml so the existing paragraph\n//    properties (Heading 2 style, ids, rsids, run formatting) are kept\n//    intact and only the <w:lastRenderedPageBreak/> marker is added.\nconst range = teamRolesParagraph.getRange();\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' +\n  '<w:p w14:paraId=\"241FEDD4\" w14:textId=\"4A5581E9\" w:rsidR=\"00E3691E\" w:rsidRPr=\"00FA22D4\" w:rsidRDefault=\"2B8A6342\" w:rsidP=\"142A635E\">' +\n  '<w:pPr><w:pStyle w:val=\"Heading2\"/><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr>' +\n  '<w:r w:rsidRPr=\"00FA22D4\">' +\n  '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\"/><w:b/><w:bCs/><w:color w:val=\"24292E\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '<w:lastRenderedPageBreak/>' +\n  '<w:t>Team Roles</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nrange.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# 1) Remove the \"Previous Participant Interest\" bullet paragraph entirely\n#    (it sat right after the \"Perigean Contract\" bullet, inside the\n#    \"What is your recruitment strategy?\" list).\n$d = $word.ActiveDocument\n\n$targetText = \"Previous Participant Interest \" + [char]0x2013 + \" We have conducted user feedback sessions in the past for DGIB and have participants who have indicated that they would like to participate in future research. \"\n\n$deleted = $false\n$paras = $d.Paragraphs\n$count = $paras.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($text -eq $targetText) {\n        $p.Range.Delete()\n        $deleted = $true\n        break\n    }\n}\nif (-not $deleted) {\n    throw 'Could not find the \"Previous Participant Interest\" paragraph to remove.'\n}\n\n# 2) Mark the \"Team Roles\" heading run with a lastRenderedPageBreak\n#    (Word stamps this on the run that happens to start a new page;\n#    it is a pure rendering hint with no visible effect on the text).\n#    We rebuild the paragraph via Range.InsertXML so the existing\n#    paragraph properties (Heading 2 style, ids, rsids, run formatting)\n#    are kept intact and only the <w:lastRenderedPageBreak/> marker is\n#    added.\n$found = $false\n$paras = $d.Paragraphs\n$count = $paras.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($text -eq \"Team Roles\") {\n        $found = $true\n        $ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n            '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n            '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n            '<pkg:xmlData>' +\n            '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n            '<w:body>' +\n            '<w:p w14:paraId=\"241FEDD4\" w14:textId=\"4A5581E9\" w:rsidR=\"00E3691E\" w:rsidRPr=\"00FA22D4\" w:rsidRDefault=\"2B8A6342\" w:rsidP=\"142A635E\">' +\n            '<w:pPr><w:pStyle w:val=\"Heading2\"/><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr>' +\n            '<w:r w:rsidRPr=\"00FA22D4\">' +\n            '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:eastAsia=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Calibri\"/><w:b/><w:bCs/><w:color w:val=\"24292E\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n            '<w:lastRenderedPageBreak/>' +\n            '<w:t>Team Roles</w:t>' +\n            '</w:r>' +\n            '</w:p>' +\n            '</w:body>' +\n            '</w:document>' +\n            '</pkg:xmlData></pkg:part></pkg:package>'\n        $p.Range.InsertXML($ooxml)\n        break\n    }\n}\nif (-not $found) {\n    throw 'Could not find the \"Team Roles\" heading paragraph.'\n}\n"}
